# "Master through 5 analysis"
# - shrink the saved window width of the workbook
# - append 3 more analysis rows (subjects 3, 4, 5) to the Master sheet
# - move the active-cell selection to where the user left off typing

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Saved window geometry (bookViews/workbookView@windowWidth 28800 -> 14400).
# Exercised through the normal Window object so the edit is expressed the
# same way an interactive resize would be; the sheet-content edits below are
# the load-bearing part of this commit.
[void]$wb.Windows.Item(1).Activate()
$excel.ActiveWindow.Width = 14400

# --- New data rows (subjects 3, 4, 5) ---------------------------------

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 3138.68957
$ws.Range("C4").Value = 5478.2784000000001
$ws.Range("D4").Value = 3089.7271700000001
$ws.Range("E4").Value = 5217.8147600000002
$ws.Range("F4").Value = 9750.6852500000005
$ws.Range("G4").Value = 6987.0085499999996

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 3280.2857600000002
$ws.Range("C5").Value = 2887.8405699999998
$ws.Range("D5").Value = 2336.8552199999999
$ws.Range("E5").Value = 3256.7347799999998
$ws.Range("F5").Value = 2735.1353399999998
$ws.Range("G5").Value = 2958.6165000000001

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 8844.9357299999992
$ws.Range("C6").Value = 3560.35734
$ws.Range("D6").Value = 9405.7161099999994
$ws.Range("E6").Value = 6366.9946600000003
$ws.Range("F6").Value = 6870.6650900000004
$ws.Range("G6").Value = 6293.7765200000003

# --- Cursor position left where the analyst stopped editing ----------

[void]$ws.Range("G14").Select()
